$d = $word.ActiveDocument

# --- Page margins (section properties) ---
$ps = $d.PageSetup
$ps.LeftMargin = 76.55    # 1531 twips
$ps.RightMargin = 76.55   # 1531 twips
$ps.BottomMargin = 96.4   # 1928 twips
# TopMargin stays at 1531 twips (76.55 pt) -- unchanged

# --- "Title" style (Style15) ---
$titleStyle = $d.Styles("Title")
$titleStyle.ParagraphFormat.SpaceAfter = 0
$titleStyle.ParagraphFormat.Alignment = 0
$titleStyle.Font.Bold = $false

# --- "Subtitle" style (Style16) ---
$subtitleStyle = $d.Styles("Subtitle")
$subtitleStyle.ParagraphFormat.LineSpacingRule = 0
$subtitleStyle.ParagraphFormat.SpaceAfter = 14.15
$subtitleStyle.Font.Bold = $true

Write-Output "done"
